# Update the "dSF" (column F) values for the rows that were repulled.
# (row number -> new value), matching the worksheet's own row numbering.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    5  = -3
    7  = -2
    12 = -3
    14 = 1
    16 = -4
    17 = -8
    18 = -1
    20 = -3
    22 = 0
    27 = -5
    30 = -7
    31 = 5
    33 = -5
    35 = -2
    36 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
